$wb = $excel.ActiveWorkbook

# ---- SongData: add rows 12-51 (ConcertID 3..10), mirroring rows 2-11 ----
$songWs = $wb.Worksheets.Item("SongData")
$songWs.Cells.Item(12,1).Value = 11
$songWs.Cells.Item(12,2).Value = 'Kell egy tank'
$songWs.Cells.Item(12,3).Value = 40
$songWs.Cells.Item(12,4).Value = 0
$songWs.Cells.Item(12,5).Value = 20
$songWs.Cells.Item(12,6).Value = $false
$songWs.Cells.Item(12,7).Value = 3
$songWs.Cells.Item(13,1).Value = 12
$songWs.Cells.Item(13,2).Value = 'Kimegyek a Szigetre'
$songWs.Cells.Item(13,3).Value = 50
$songWs.Cells.Item(13,4).Value = 0
$songWs.Cells.Item(13,5).Value = 30
$songWs.Cells.Item(13,6).Value = $false
$songWs.Cells.Item(13,7).Value = 3
$songWs.Cells.Item(14,1).Value = 13
$songWs.Cells.Item(14,2).Value = 'Ejnye Ottó'
$songWs.Cells.Item(14,3).Value = 60
$songWs.Cells.Item(14,4).Value = 0
$songWs.Cells.Item(14,5).Value = 40
$songWs.Cells.Item(14,6).Value = $false
$songWs.Cells.Item(14,7).Value = 3
$songWs.Cells.Item(15,1).Value = 14
$songWs.Cells.Item(15,2).Value = 'Tahó család'
$songWs.Cells.Item(15,3).Value = 70
$songWs.Cells.Item(15,4).Value = 0
$songWs.Cells.Item(15,5).Value = 50
$songWs.Cells.Item(15,6).Value = $false
$songWs.Cells.Item(15,7).Value = 3
$songWs.Cells.Item(16,1).Value = 15
$songWs.Cells.Item(16,2).Value = 'Szegény ember kézzel nőz'
$songWs.Cells.Item(16,3).Value = 80
$songWs.Cells.Item(16,4).Value = 20
$songWs.Cells.Item(16,5).Value = 0
$songWs.Cells.Item(16,6).Value = $true
$songWs.Cells.Item(16,7).Value = 3
$songWs.Cells.Item(17,1).Value = 16
$songWs.Cells.Item(17,2).Value = 'Akkor inkább úthenger'
$songWs.Cells.Item(17,3).Value = 60
$songWs.Cells.Item(17,4).Value = 0
$songWs.Cells.Item(17,5).Value = 25
$songWs.Cells.Item(17,6).Value = $false
$songWs.Cells.Item(17,7).Value = 4
$songWs.Cells.Item(18,1).Value = 17
$songWs.Cells.Item(18,2).Value = 'Legyünk mi is pöcsfejek'
$songWs.Cells.Item(18,3).Value = 70
$songWs.Cells.Item(18,4).Value = 0
$songWs.Cells.Item(18,5).Value = 35
$songWs.Cells.Item(18,6).Value = $false
$songWs.Cells.Item(18,7).Value = 4
$songWs.Cells.Item(19,1).Value = 18
$songWs.Cells.Item(19,2).Value = 'Télapó itt van'
$songWs.Cells.Item(19,3).Value = 80
$songWs.Cells.Item(19,4).Value = 0
$songWs.Cells.Item(19,5).Value = 45
$songWs.Cells.Item(19,6).Value = $false
$songWs.Cells.Item(19,7).Value = 4
$songWs.Cells.Item(20,1).Value = 19
$songWs.Cells.Item(20,2).Value = 'Kiskunszicíla'
$songWs.Cells.Item(20,3).Value = 90
$songWs.Cells.Item(20,4).Value = 0
$songWs.Cells.Item(20,5).Value = 55
$songWs.Cells.Item(20,6).Value = $false
$songWs.Cells.Item(20,7).Value = 4
$songWs.Cells.Item(21,1).Value = 20
$songWs.Cells.Item(21,2).Value = 'Tábor az erdő mélyén'
$songWs.Cells.Item(21,3).Value = 100
$songWs.Cells.Item(21,4).Value = 20
$songWs.Cells.Item(21,5).Value = 0
$songWs.Cells.Item(21,6).Value = $true
$songWs.Cells.Item(21,7).Value = 4
$songWs.Cells.Item(22,1).Value = 21
$songWs.Cells.Item(22,2).Value = 'Kell egy tank'
$songWs.Cells.Item(22,3).Value = 40
$songWs.Cells.Item(22,4).Value = 0
$songWs.Cells.Item(22,5).Value = 20
$songWs.Cells.Item(22,6).Value = $false
$songWs.Cells.Item(22,7).Value = 5
$songWs.Cells.Item(23,1).Value = 22
$songWs.Cells.Item(23,2).Value = 'Kimegyek a Szigetre'
$songWs.Cells.Item(23,3).Value = 50
$songWs.Cells.Item(23,4).Value = 0
$songWs.Cells.Item(23,5).Value = 30
$songWs.Cells.Item(23,6).Value = $false
$songWs.Cells.Item(23,7).Value = 5
$songWs.Cells.Item(24,1).Value = 23
$songWs.Cells.Item(24,2).Value = 'Ejnye Ottó'
$songWs.Cells.Item(24,3).Value = 60
$songWs.Cells.Item(24,4).Value = 0
$songWs.Cells.Item(24,5).Value = 40
$songWs.Cells.Item(24,6).Value = $false
$songWs.Cells.Item(24,7).Value = 5
$songWs.Cells.Item(25,1).Value = 24
$songWs.Cells.Item(25,2).Value = 'Tahó család'
$songWs.Cells.Item(25,3).Value = 70
$songWs.Cells.Item(25,4).Value = 0
$songWs.Cells.Item(25,5).Value = 50
$songWs.Cells.Item(25,6).Value = $false
$songWs.Cells.Item(25,7).Value = 5
$songWs.Cells.Item(26,1).Value = 25
$songWs.Cells.Item(26,2).Value = 'Szegény ember kézzel nőz'
$songWs.Cells.Item(26,3).Value = 80
$songWs.Cells.Item(26,4).Value = 20
$songWs.Cells.Item(26,5).Value = 0
$songWs.Cells.Item(26,6).Value = $true
$songWs.Cells.Item(26,7).Value = 5
$songWs.Cells.Item(27,1).Value = 26
$songWs.Cells.Item(27,2).Value = 'Akkor inkább úthenger'
$songWs.Cells.Item(27,3).Value = 60
$songWs.Cells.Item(27,4).Value = 0
$songWs.Cells.Item(27,5).Value = 25
$songWs.Cells.Item(27,6).Value = $false
$songWs.Cells.Item(27,7).Value = 6
$songWs.Cells.Item(28,1).Value = 27
$songWs.Cells.Item(28,2).Value = 'Legyünk mi is pöcsfejek'
$songWs.Cells.Item(28,3).Value = 70
$songWs.Cells.Item(28,4).Value = 0
$songWs.Cells.Item(28,5).Value = 35
$songWs.Cells.Item(28,6).Value = $false
$songWs.Cells.Item(28,7).Value = 6
$songWs.Cells.Item(29,1).Value = 28
$songWs.Cells.Item(29,2).Value = 'Télapó itt van'
$songWs.Cells.Item(29,3).Value = 80
$songWs.Cells.Item(29,4).Value = 0
$songWs.Cells.Item(29,5).Value = 45
$songWs.Cells.Item(29,6).Value = $false
$songWs.Cells.Item(29,7).Value = 6
$songWs.Cells.Item(30,1).Value = 29
$songWs.Cells.Item(30,2).Value = 'Kiskunszicíla'
$songWs.Cells.Item(30,3).Value = 90
$songWs.Cells.Item(30,4).Value = 0
$songWs.Cells.Item(30,5).Value = 55
$songWs.Cells.Item(30,6).Value = $false
$songWs.Cells.Item(30,7).Value = 6
$songWs.Cells.Item(31,1).Value = 30
$songWs.Cells.Item(31,2).Value = 'Tábor az erdő mélyén'
$songWs.Cells.Item(31,3).Value = 100
$songWs.Cells.Item(31,4).Value = 20
$songWs.Cells.Item(31,5).Value = 0
$songWs.Cells.Item(31,6).Value = $true
$songWs.Cells.Item(31,7).Value = 6
$songWs.Cells.Item(32,1).Value = 31
$songWs.Cells.Item(32,2).Value = 'Kell egy tank'
$songWs.Cells.Item(32,3).Value = 40
$songWs.Cells.Item(32,4).Value = 0
$songWs.Cells.Item(32,5).Value = 20
$songWs.Cells.Item(32,6).Value = $false
$songWs.Cells.Item(32,7).Value = 7
$songWs.Cells.Item(33,1).Value = 32
$songWs.Cells.Item(33,2).Value = 'Kimegyek a Szigetre'
$songWs.Cells.Item(33,3).Value = 50
$songWs.Cells.Item(33,4).Value = 0
$songWs.Cells.Item(33,5).Value = 30
$songWs.Cells.Item(33,6).Value = $false
$songWs.Cells.Item(33,7).Value = 7
$songWs.Cells.Item(34,1).Value = 33
$songWs.Cells.Item(34,2).Value = 'Ejnye Ottó'
$songWs.Cells.Item(34,3).Value = 60
$songWs.Cells.Item(34,4).Value = 0
$songWs.Cells.Item(34,5).Value = 40
$songWs.Cells.Item(34,6).Value = $false
$songWs.Cells.Item(34,7).Value = 7
$songWs.Cells.Item(35,1).Value = 34
$songWs.Cells.Item(35,2).Value = 'Tahó család'
$songWs.Cells.Item(35,3).Value = 70
$songWs.Cells.Item(35,4).Value = 0
$songWs.Cells.Item(35,5).Value = 50
$songWs.Cells.Item(35,6).Value = $false
$songWs.Cells.Item(35,7).Value = 7
$songWs.Cells.Item(36,1).Value = 35
$songWs.Cells.Item(36,2).Value = 'Szegény ember kézzel nőz'
$songWs.Cells.Item(36,3).Value = 80
$songWs.Cells.Item(36,4).Value = 20
$songWs.Cells.Item(36,5).Value = 0
$songWs.Cells.Item(36,6).Value = $true
$songWs.Cells.Item(36,7).Value = 7
$songWs.Cells.Item(37,1).Value = 36
$songWs.Cells.Item(37,2).Value = 'Akkor inkább úthenger'
$songWs.Cells.Item(37,3).Value = 60
$songWs.Cells.Item(37,4).Value = 0
$songWs.Cells.Item(37,5).Value = 25
$songWs.Cells.Item(37,6).Value = $false
$songWs.Cells.Item(37,7).Value = 8
$songWs.Cells.Item(38,1).Value = 37
$songWs.Cells.Item(38,2).Value = 'Legyünk mi is pöcsfejek'
$songWs.Cells.Item(38,3).Value = 70
$songWs.Cells.Item(38,4).Value = 0
$songWs.Cells.Item(38,5).Value = 35
$songWs.Cells.Item(38,6).Value = $false
$songWs.Cells.Item(38,7).Value = 8
$songWs.Cells.Item(39,1).Value = 38
$songWs.Cells.Item(39,2).Value = 'Télapó itt van'
$songWs.Cells.Item(39,3).Value = 80
$songWs.Cells.Item(39,4).Value = 0
$songWs.Cells.Item(39,5).Value = 45
$songWs.Cells.Item(39,6).Value = $false
$songWs.Cells.Item(39,7).Value = 8
$songWs.Cells.Item(40,1).Value = 39
$songWs.Cells.Item(40,2).Value = 'Kiskunszicíla'
$songWs.Cells.Item(40,3).Value = 90
$songWs.Cells.Item(40,4).Value = 0
$songWs.Cells.Item(40,5).Value = 55
$songWs.Cells.Item(40,6).Value = $false
$songWs.Cells.Item(40,7).Value = 8
$songWs.Cells.Item(41,1).Value = 40
$songWs.Cells.Item(41,2).Value = 'Tábor az erdő mélyén'
$songWs.Cells.Item(41,3).Value = 100
$songWs.Cells.Item(41,4).Value = 20
$songWs.Cells.Item(41,5).Value = 0
$songWs.Cells.Item(41,6).Value = $true
$songWs.Cells.Item(41,7).Value = 8
$songWs.Cells.Item(42,1).Value = 41
$songWs.Cells.Item(42,2).Value = 'Kell egy tank'
$songWs.Cells.Item(42,3).Value = 40
$songWs.Cells.Item(42,4).Value = 0
$songWs.Cells.Item(42,5).Value = 20
$songWs.Cells.Item(42,6).Value = $false
$songWs.Cells.Item(42,7).Value = 9
$songWs.Cells.Item(43,1).Value = 42
$songWs.Cells.Item(43,2).Value = 'Kimegyek a Szigetre'
$songWs.Cells.Item(43,3).Value = 50
$songWs.Cells.Item(43,4).Value = 0
$songWs.Cells.Item(43,5).Value = 30
$songWs.Cells.Item(43,6).Value = $false
$songWs.Cells.Item(43,7).Value = 9
$songWs.Cells.Item(44,1).Value = 43
$songWs.Cells.Item(44,2).Value = 'Ejnye Ottó'
$songWs.Cells.Item(44,3).Value = 60
$songWs.Cells.Item(44,4).Value = 0
$songWs.Cells.Item(44,5).Value = 40
$songWs.Cells.Item(44,6).Value = $false
$songWs.Cells.Item(44,7).Value = 9
$songWs.Cells.Item(45,1).Value = 44
$songWs.Cells.Item(45,2).Value = 'Tahó család'
$songWs.Cells.Item(45,3).Value = 70
$songWs.Cells.Item(45,4).Value = 0
$songWs.Cells.Item(45,5).Value = 50
$songWs.Cells.Item(45,6).Value = $false
$songWs.Cells.Item(45,7).Value = 9
$songWs.Cells.Item(46,1).Value = 45
$songWs.Cells.Item(46,2).Value = 'Szegény ember kézzel nőz'
$songWs.Cells.Item(46,3).Value = 80
$songWs.Cells.Item(46,4).Value = 20
$songWs.Cells.Item(46,5).Value = 0
$songWs.Cells.Item(46,6).Value = $true
$songWs.Cells.Item(46,7).Value = 9
$songWs.Cells.Item(47,1).Value = 46
$songWs.Cells.Item(47,2).Value = 'Akkor inkább úthenger'
$songWs.Cells.Item(47,3).Value = 60
$songWs.Cells.Item(47,4).Value = 0
$songWs.Cells.Item(47,5).Value = 25
$songWs.Cells.Item(47,6).Value = $false
$songWs.Cells.Item(47,7).Value = 10
$songWs.Cells.Item(48,1).Value = 47
$songWs.Cells.Item(48,2).Value = 'Legyünk mi is pöcsfejek'
$songWs.Cells.Item(48,3).Value = 70
$songWs.Cells.Item(48,4).Value = 0
$songWs.Cells.Item(48,5).Value = 35
$songWs.Cells.Item(48,6).Value = $false
$songWs.Cells.Item(48,7).Value = 10
$songWs.Cells.Item(49,1).Value = 48
$songWs.Cells.Item(49,2).Value = 'Télapó itt van'
$songWs.Cells.Item(49,3).Value = 80
$songWs.Cells.Item(49,4).Value = 0
$songWs.Cells.Item(49,5).Value = 45
$songWs.Cells.Item(49,6).Value = $false
$songWs.Cells.Item(49,7).Value = 10
$songWs.Cells.Item(50,1).Value = 49
$songWs.Cells.Item(50,2).Value = 'Kiskunszicíla'
$songWs.Cells.Item(50,3).Value = 90
$songWs.Cells.Item(50,4).Value = 0
$songWs.Cells.Item(50,5).Value = 55
$songWs.Cells.Item(50,6).Value = $false
$songWs.Cells.Item(50,7).Value = 10
$songWs.Cells.Item(51,1).Value = 50
$songWs.Cells.Item(51,2).Value = 'Tábor az erdő mélyén'
$songWs.Cells.Item(51,3).Value = 100
$songWs.Cells.Item(51,4).Value = 20
$songWs.Cells.Item(51,5).Value = 0
$songWs.Cells.Item(51,6).Value = $true
$songWs.Cells.Item(51,7).Value = 10

# copy formatting for the new data rows from the last existing row (row 11)
$songWs.Range("A11:G11").Copy()
$songWs.Range("A12:G51").PasteSpecial(-4122)

# rows 52-101: blank placeholder rows (A:F only), matching the sheet's style
$songWs.Range("A2:F2").Copy()
$songWs.Range("A52:F101").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- ConcertData: add rows 4-12 for the new Falunap 2..10 concerts ----
$concertWs = $wb.Worksheets.Item("ConcertData")
$concertWs.Cells.Item(4,1).Value = 3
$concertWs.Cells.Item(4,2).Value = 'Falunap 2'
$concertWs.Cells.Item(4,3).Value = 90
$concertWs.Cells.Item(4,4).Value = 200
$concertWs.Cells.Item(5,1).Value = 4
$concertWs.Cells.Item(5,2).Value = 'Falunap 3'
$concertWs.Cells.Item(5,3).Value = 110
$concertWs.Cells.Item(5,4).Value = 250
$concertWs.Cells.Item(6,1).Value = 5
$concertWs.Cells.Item(6,2).Value = 'Falunap 4'
$concertWs.Cells.Item(6,3).Value = 130
$concertWs.Cells.Item(6,4).Value = 300
$concertWs.Cells.Item(7,1).Value = 6
$concertWs.Cells.Item(7,2).Value = 'Falunap 5'
$concertWs.Cells.Item(7,3).Value = 150
$concertWs.Cells.Item(7,4).Value = 350
$concertWs.Cells.Item(8,1).Value = 7
$concertWs.Cells.Item(8,2).Value = 'Falunap 6'
$concertWs.Cells.Item(8,3).Value = 170
$concertWs.Cells.Item(8,4).Value = 400
$concertWs.Cells.Item(9,1).Value = 8
$concertWs.Cells.Item(9,2).Value = 'Falunap 7'
$concertWs.Cells.Item(9,3).Value = 190
$concertWs.Cells.Item(9,4).Value = 450
$concertWs.Cells.Item(10,1).Value = 9
$concertWs.Cells.Item(10,2).Value = 'Falunap 8'
$concertWs.Cells.Item(10,3).Value = 210
$concertWs.Cells.Item(10,4).Value = 500
$concertWs.Cells.Item(11,1).Value = 10
$concertWs.Cells.Item(11,2).Value = 'Falunap 9'
$concertWs.Cells.Item(11,3).Value = 230
$concertWs.Cells.Item(11,4).Value = 550
$concertWs.Cells.Item(12,1).Value = 11
$concertWs.Cells.Item(12,2).Value = 'Falunap 10'
$concertWs.Cells.Item(12,3).Value = 250
$concertWs.Cells.Item(12,4).Value = 600

# copy formatting for the new rows from the last existing row (row 3)
$concertWs.Range("A3:D3").Copy()
$concertWs.Range("A4:D12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- MerchData bugfix: Table & Stickers (level 1) now cost 50 to unlock ----
$merchWs = $wb.Worksheets.Item("MerchData")
$merchWs.Cells.Item(2,5).Value = 50
$merchWs.Cells.Item(5,5).Value = 50

